# ---------------------------------------------------------------------------
# Adds a new "Signal_Value_123" column (AJ) to the mounted-pipeline workbook
# and updates all downstream derived values (Step2_Sj cumulative sums and the
# Step3_DataPts_* "Point_Exceeds_Cumulative_Value" lookups) to reflect it.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step1_Data (sheet1) and Step2_Sj (sheet2): both gain a new column AJ with
# header "Signal_Value_123", formatted like the existing AI header cell.
# ---------------------------------------------------------------------------
$sheetNames = @("Step1_Data", "Step2_Sj")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("AJ1").Value = "Signal_Value_123"
    $ws.Range("AI1").Copy() | Out-Null
    $ws.Range("AJ1").PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------------
# Step1_Data: row 3 ("signal segment 2") gets renormalised intensity values
# now that Signal_Value_123 carries some of the weight, and a new AJ3 value.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Step1_Data")

$row3Updates = @{
    "D3"  = 0
    "E3"  = 0.1307621379239056
    "F3"  = 0.09285566823687409
    "G3"  = 0.1880508824833489
    "K3"  = 0.01581102621035328
    "L3"  = 0.01545534223231766
    "M3"  = 0.06981312204652913
    "N3"  = 0.003284179928965613
    "O3"  = 0.1155616491805911
    "P3"  = 0.007232141138927713
    "Q3"  = 0.02297810451142166
    "R3"  = 0.001962844361772281
    "S3"  = 0.0423201944409661
    "T3"  = 0.04046204704280144
    "V3"  = 0.0004119555421370758
    "Z3"  = 0.02054955909943972
    "AA3" = 0.05370304966475895
    "AC3" = 0.04880398628564083
    "AD3" = 0.01539321761865344
    "AE3" = 0.02074791430278316
    "AF3" = 0.05029949255128497
    "AG3" = 0.0349771104501679
    "AH3" = 0.008171355541042882
}
foreach ($addr in $row3Updates.Keys) {
    $ws1.Range($addr).Value = $row3Updates[$addr]
}

# New AJ column for Step1_Data, rows 2-11 (only row 3 is non-zero).
$sheet1AJ = @{
    2  = 0
    3  = 0.000393019205316371
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
}
foreach ($r in $sheet1AJ.Keys) {
    $ws1.Cells.Item($r, 36).Value = $sheet1AJ[$r]
}

# ---------------------------------------------------------------------------
# Step2_Sj: row 3 is the running cumulative sum of Step1_Data row 3, so every
# value from D3 onward is recomputed; a new AJ3 (= 1.0) is appended.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Step2_Sj")

$sheet2Row3 = @{
    "D3"  = 0
    "E3"  = 0.1307621379239056
    "F3"  = 0.2236178061607796
    "G3"  = 0.4116686886441286
    "H3"  = 0.4116686886441286
    "I3"  = 0.4116686886441286
    "J3"  = 0.4116686886441286
    "K3"  = 0.4274797148544819
    "L3"  = 0.4429350570867995
    "M3"  = 0.5127481791333287
    "N3"  = 0.5160323590622943
    "O3"  = 0.6315940082428854
    "P3"  = 0.638826149381813
    "Q3"  = 0.6618042538932347
    "R3"  = 0.663767098255007
    "S3"  = 0.7060872926959731
    "T3"  = 0.7465493397387746
    "U3"  = 0.7465493397387746
    "V3"  = 0.7469612952809116
    "W3"  = 0.7469612952809116
    "X3"  = 0.7469612952809116
    "Y3"  = 0.7469612952809116
    "Z3"  = 0.7675108543803513
    "AA3" = 0.8212139040451102
    "AB3" = 0.8212139040451102
    "AC3" = 0.8700178903307511
    "AD3" = 0.8854111079494045
    "AE3" = 0.9061590222521876
    "AF3" = 0.9564585148034725
    "AG3" = 0.9914356252536405
    "AH3" = 0.9996069807946834
    "AI3" = 0.9996069807946834
    "AJ3" = 0.9999999999999998
}
foreach ($addr in $sheet2Row3.Keys) {
    $ws2.Range($addr).Value = $sheet2Row3[$addr]
}

# New AJ column for Step2_Sj: rows whose Step1_Data AJ contribution is 0 just
# carry forward the same cumulative total that was already in column AI.
$sheet2AJ = @{
    2  = 1
    4  = 0.9999999999999999
    5  = 0.9999999999999998
    6  = 0.9999999999999997
    7  = 1
    8  = 0.9999999999999999
    9  = 1
    10 = 0.9999999999999999
    11 = 0.9999999999999998
}
foreach ($r in $sheet2AJ.Keys) {
    $ws2.Cells.Item($r, 36).Value = $sheet2AJ[$r]
}

# ---------------------------------------------------------------------------
# Step3_DataPts_* sheets: "Point_Exceeds_Cumulative_Value" (F3) is looked up
# from Step2_Sj's recomputed row 3, so each threshold sheet's F3 is refreshed.
# ---------------------------------------------------------------------------
$step3Updates = @{
    "Step3_DataPts_0.5" = 0.5127481791333287
    "Step3_DataPts_0.7" = 0.7060872926959731
    "Step3_DataPts_0.8" = 0.8212139040451102
    "Step3_DataPts_0.9" = 0.9061590222521876
}
foreach ($name in $step3Updates.Keys) {
    $ws3 = $wb.Worksheets.Item($name)
    $ws3.Range("F3").Value = $step3Updates[$name]
}
